$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header pin references in rows 2-3 (A/D columns):
#    A2: 5V -> PIN 4 ; A3: GND -> PIN 5
#    D2: 5V -> PIN 6 ; D3: GND -> PIN 7
$ws.Range("A2").Value = "PIN 4"
$ws.Range("A3").Value = "PIN 5"
$ws.Range("D2").Value = "PIN 6"
$ws.Range("D3").Value = "PIN 7"

# 2. Remove the whole "MOTOR DC 1 DAN 2 / MOTOR DC 3 DAN 4" block (old rows 5:7).
#    Deleting shifts everything below up by 3 rows, which naturally lines the
#    remaining blocks up with their new target row numbers.
$ws.Rows("5:7").Delete()

# 3. Materialize the blank spacer row at row 4 (previously an untouched gap)
#    with an explicit (borderless) style, matching the new layout.
$ws.Range("A4:E4").Borders.LineStyle = -4142

# 4. Update the selection shown when the sheet is reopened.
$ws.Range("K19").Select()
